# Append the new resale-number observation row (2025-02-11 22:03:14) as row 64.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 64

# Columns A ("2025-02-11") and D ("06") look like a date / number respectively,
# so a leading apostrophe is used to force Excel to store them as literal text
# (same effect as pre-formatting the cell as Text before typing the value).
$ws.Cells.Item($row, 1).Value  = "'2025-02-11"
$ws.Cells.Item($row, 2).Value  = "22:03:14"
$ws.Cells.Item($row, 3).Value  = "Tuesday"
$ws.Cells.Item($row, 4).Value  = "'06"
$ws.Cells.Item($row, 5).Value  = 127946
$ws.Cells.Item($row, 6).Value  = 142006
$ws.Cells.Item($row, 7).Value  = 169319
$ws.Cells.Item($row, 8).Value  = 158826
$ws.Cells.Item($row, 9).Value  = -1
$ws.Cells.Item($row, 10).Value = 144488
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191660
$ws.Cells.Item($row, 14).Value = 115021
$ws.Cells.Item($row, 15).Value = 44984
$ws.Cells.Item($row, 16).Value = 28544
$ws.Cells.Item($row, 17).Value = 65251
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 43992
$ws.Cells.Item($row, 20).Value = -1
